$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.092.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.655.42'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5237'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.65%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06346'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.40'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07805'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('E12').Value = '  +1.58%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.637.95'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.34%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5480'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.93%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0₅8255'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.88%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.50'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.117.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('E18').Value = '  -0.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.582'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.91'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.07'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.033'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '142.08'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1234'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.93%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.242'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.07'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.22%  '
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05876'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('E30').Value = '  +0.34%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.537'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.268'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.582'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.87%  '
$ws.Range('E34').Value = '  -1.44%  '
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.780'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.5722'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.13%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01616'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.783'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.8459'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.27%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '103.23'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.027.53'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.800.16'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('E46').Value = '  -0.44%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4316'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05151'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.468'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.802'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.09666'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.35%  '
